$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MCT-1A-Gestão integrada"
$ws.Range("C3").Value = "MCT-1A-Gestão integrada"
$ws.Range("E6").Value = "-"
$ws.Range("E7").Value = "-"
